$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement dataset for rows 2-31 (A:H) reflecting the "May 9th" data update:
# 9 new samples inserted at the top (after the header), the previously-existing 20
# samples shifted down, and 1 new sample appended at the end.
$rows = @(
  @{ A=0; B="falling"; C=-2.822264432907104; D=3.932364892959594; E=0.4758871763944626; F=0.0328340083360672; G=0.0253509078174829; H=0.0371100641787052 },
  @{ A=100; B="falling"; C=-2.760175883769989; D=3.927277326583862; E=0.4960161261260509; F=0.0198531206697225; G=0.0163406450301408; H=-0.009468411095440299 },
  @{ A=200; B="falling"; C=-2.64756965637207; D=3.941246557235718; E=0.6029229372739793; F=-0.0001527163112768; G=0.0464257597923278; H=0.0154243474826216 },
  @{ A=300; B="falling"; C=-2.628955054283141; D=3.869387286901474; E=0.6939661800861359; F=0.0108428578823804; G=-0.015118914656341; H=0.1504255682229995 },
  @{ A=400; B="falling"; C=-2.694027137756347; D=3.838324213027954; E=0.7007610917091369; F=0.0204639863222837; G=-0.0316122770309448; H=0.0610865242779254 },
  @{ A=500; B="falling"; C=-2.807751727104188; D=3.847099477052689; E=0.6076438263058663; F=-0.009010262787342; G=-0.0128281703218817; H=0.0167987942695617 },
  @{ A=600; B="falling"; C=-2.889829158782959; D=3.764959990978241; E=0.4738430827856064; F=-0.0468839071691036; G=0.052381694316864; H=0.0694859251379966 },
  @{ A=700; B="falling"; C=-3.067947173118591; D=3.764124345779419; E=0.4743617072701454; F=0.0027488935738801; G=0.09178250283002851; H=0.0717766657471656 },
  @{ A=800; B="falling"; C=-3.113880395889281; D=3.648524475097656; E=0.3546387374401094; F=0.0128281703218817; G=0.0387899428606033; H=0.0054977871477603 },
  @{ A=900; B="falling"; C=-2.830316853523255; D=3.761647629737854; E=0.4785371914505959; F=0.1214094683527946; G=-0.0080939643085002; H=0.0326812900602817 },
  @{ A=1000; B="falling"; C=-2.873888826370239; D=3.739926481246949; E=0.2519635170698165; F=-0.061391957104206; G=-0.0507018156349659; H=-0.0813977941870689 },
  @{ A=1100; B="falling"; C=-3.17759734392166; D=3.094935894012452; E=-1.055553257465361; F=0.0126754539087414; G=-0.266184538602829; H=0.0452040284872055 },
  @{ A=1200; B="falling"; C=-3.781754684448242; D=2.629349136352539; E=-1.783702611923218; F=0.0755945742130279; G=-0.1324050426483154; H=0.1108720451593399 },
  @{ A=1300; B="falling"; C=-3.911327278614044; D=1.632070899009704; E=-2.105818438529968; F=0.2214386463165283; G=-0.4340197443962097; H=-0.0148134818300604 },
  @{ A=1400; B="falling"; C=-5.623295831680297; D=0.7180684566497808; E=-2.688220548629761; F=0.2553416788578033; G=-0.5070181488990784; H=0.0436768643558025 },
  @{ A=1500; B="falling"; C=-6.422207009792328; D=0.5463402450084687; E=-3.518551957607269; F=0.8651378750801086; G=0.230448916554451; H=-0.6892086863517761 },
  @{ A=1600; B="falling"; C=-5.178159475326543; D=1.582323312759396; E=-3.006135582923893; F=1.613447785377502; G=-0.0974330082535743; H=-2.265546560287476 },
  @{ A=1700; B="falling"; C=-3.926182019710541; D=3.611184996366497; E=0.6073530614375975; F=0.1902845203876495; G=1.200044751167297; H=-3.24567985534668 },
  @{ A=1800; B="falling"; C=-4.239339423179628; D=4.14080636501312; E=3.278265589475629; F=-0.5577199459075928; G=2.010357618331909; H=0.8011497855186462 },
  @{ A=1900; B="falling"; C=-1.282989490032158; D=2.413680851459504; E=0.6543272763490596; F=-3.031724214553833; G=-0.6270532011985779; H=-4.525289535522461 },
  @{ A=2000; B="falling"; C=7.509865951538076; D=2.822640895843509; E=-1.422052669525146; F=-2.38100004196167; G=0.087353728711605; H=-1.472948789596558 },
  @{ A=2100; B="falling"; C=2.871775150299069; D=4.098878204822539; E=-1.140487685799597; F=0.8362745046615601; G=2.635119915008545; H=0.4735732674598694 },
  @{ A=2200; B="falling"; C=3.782487344741822; D=4.224413537979126; E=-1.672759181261064; F=-0.3104722499847412; G=-0.0462730415165424; H=-0.327729195356369 },
  @{ A=2300; B="falling"; C=-0.05776283740998212; D=3.633949923515318; E=0.03897095024586128; F=-0.0215329993516206; G=0.9650143980979921; H=0.0117591563612222 },
  @{ A=2400; B="falling"; C=3.522613048553472; D=5.809941577911379; E=-2.533208680152897; F=-0.5355761051177979; G=-0.6331618428230286; H=-0.1020144969224929 },
  @{ A=2500; B="falling"; C=4.382918024063111; D=3.592077732086182; E=-2.185777962207795; F=0.0649044290184974; G=0.545655369758606; H=0.1579086631536483 },
  @{ A=2600; B="falling"; C=2.164362668991095; D=3.682830810546873; E=-3.28418397903442; F=0.1128573566675186; G=0.2232712507247924; H=0.1513418704271316 },
  @{ A=2700; B="falling"; C=1.118086504936218; D=3.899120330810547; E=-2.249887198209762; F=0.1999056488275528; G=-0.3659082949161529; H=-0.0462730415165424 },
  @{ A=2800; B="falling"; C=0.6892168045043942; D=4.255544376373291; E=-1.803213787078858; F=0.2883284091949463; G=-0.233655959367752; H=-0.078343465924263 },
  @{ A=2900; B="falling"; C=0.8177140951156612; D=4.59937185049057; E=-1.870238688588142; F=0.3949243724346161; G=0.0597120784223079; H=0.08445212244987479 }
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value2 = $row.A
  $ws.Cells.Item($r, 2).Value2 = $row.B
  $ws.Cells.Item($r, 3).Value2 = $row.C
  $ws.Cells.Item($r, 4).Value2 = $row.D
  $ws.Cells.Item($r, 5).Value2 = $row.E
  $ws.Cells.Item($r, 6).Value2 = $row.F
  $ws.Cells.Item($r, 7).Value2 = $row.G
  $ws.Cells.Item($r, 8).Value2 = $row.H
  $r++
}

Write-Output "Updated $($rows.Count) rows; last row = $($r - 1)"